# Applies the "Added new deployment file" change:
#   - source filename for the File -> RAW checks (CIFF_*, rows 2-8, column C)
#     changes from contact_info_20240709.csv to contact_info_20240913.csv
#   - the RAW table (ETL_AUTO.CONTACT_INFO.CONTACT_INFO_RAW) used as the
#     target of the File -> RAW checks (rows 2-8, column H) and as the
#     source of the RAW -> BRONZE checks (rows 9-16, column C) is renamed
#     to SAMPLEDB.CONTACT_INFO.CONTACT_INFO_RAW
#   - the BRONZE table (ETL_AUTO.CONTACT_INFO.CONTACT_INFO_BRONZE) used as
#     the target of the RAW -> BRONZE checks (rows 9-16, column H) is
#     renamed to SAMPLEDB.CONTACT_INFO.CONTACT_INFO_BRONZE

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data_validation")

$newRaw = "SAMPLEDB.CONTACT_INFO.CONTACT_INFO_RAW"
$newBronze = "SAMPLEDB.CONTACT_INFO.CONTACT_INFO_BRONZE"
$newFile = "contact_info_20240913.csv"

for ($r = 9; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = $newRaw
    $ws.Cells.Item($r, 8).Value = $newBronze
}

for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = $newFile
    $ws.Cells.Item($r, 8).Value = $newRaw
}

# Refresh the view: move the active selection onto the new source column
# (the frozen pane's top-left cell follows the selection automatically).
$ws.Activate()
$ws.Range("C3:C8").Select()
